# Applies the 8/1/2014 updates to the "Integrated Register" risk register:
#  - fixes a typo in the quality-monitoring-tools risk description (row 7)
#  - closes out the "Users are not able to perform their job" risk (row 9)
#  - adds a brand-new risk row (row 11) about engineers supporting multiple
#    projects, with its mitigation plan, dates, probability/impact, etc.
#  - tidies up related formatting (row heights, Probability column wrap,
#    and the new row's Context cell alignment) and moves the active
#    selection to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the typo "reuire" -> "require" in row 7's risk description.
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "Changes, such as changes to quality monitoring tools may impact functionality and require engineering changes"

# ---------------------------------------------------------------------
# 2. Close out row 9's risk: mark Status "Closed", record the mitigation
#    stop date, and append the closure note to the Threshold/Trigger cell.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Closed"
$ws.Range("R9").Value = (Get-Date -Year 2014 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("T9").Value = "Users are not able to perform their job`nRisk is Closed 7/30/14"

# ---------------------------------------------------------------------
# 3. Populate the new risk in row 11.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "open"
$ws.Range("C11").Value = "Engineers supporting multiple projects will have negative impact on delivery of eCL changes"
$ws.Range("D11").Value = (Get-Date -Year 2014 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E11").Value = "Resources are supporting eCL, Performance Scorecard, A&E, IQS, Aspect"
$ws.Range("F11").Value = "Tim"
$ws.Range("G11").Value = "Tim"
$ws.Range("H11").Value = "eCL team"
$ws.Range("I11").Value = "System"
$ws.Range("J11").Value = 0.85
$ws.Range("K11").Value = 5
$ws.Range("P11").Value = "Balance work load across engineers to the greatest extent possible to prevent distractions."
$ws.Range("Q11").Value = (Get-Date -Year 2014 -Month 7 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("S11").Value = "Balance work load across engineers to the greatest extent possible to prevent distractions."

# Re-assert the Risk Exposure Ranking formula so it recalculates now that
# Probability/Impact are real numbers instead of the blank-row placeholder.
$ws.Range("M11").Formula = $ws.Range("M11").Formula

# Context cell for the new row reads left-aligned (rather than the table's
# default general alignment) while keeping top vertical alignment + wrap.
$ws.Range("E11").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("E11").VerticalAlignment = -4160     # xlVAlignTop
$ws.Range("E11").WrapText = $true

# ---------------------------------------------------------------------
# 4. Formatting touch-ups that follow from the above edits.
# ---------------------------------------------------------------------
# The Probability column no longer force-wraps its percentages.
$ws.Range("J4:J37").WrapText = $false

# Row heights grow to fit the longer wrapped text.
$ws.Rows.Item(7).RowHeight = 63
$ws.Rows.Item(9).RowHeight = 88.2
$ws.Rows.Item(11).RowHeight = 63

# ---------------------------------------------------------------------
# 5. Move the active selection to the new row, as the author left it.
# ---------------------------------------------------------------------
$ws.Range("A11").Select()
